# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (rows 16-24, columns C:F) is being re-sorted by
# "Periodo Mora" (column E) - the 1710 period moves to the top, followed by
# the four 1802 rows, then the four 1803 rows - while keeping the same
# employee (N Doc Trabajador / Nombre Trabajador) and Valor Mora (F) tied to
# its own row. Columns B (Tipo Doc Trabajador = "CC") and G (Salario Basico
# = 781242) are identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; Doc = "9148752";  Nombre = "HUMBERTO DE JESUS MULET HERNANDEZ";   Periodo = "1710"; Mora = 22624 },
    @{ Row = 17; Doc = "22800539"; Nombre = "DARLY YANETH PEREA TUIRAN";           Periodo = "1802"; Mora = 31249 },
    @{ Row = 18; Doc = "73099457"; Nombre = "RAFAEL SALAZAR MACHADO";             Periodo = "1802"; Mora = 31249 },
    @{ Row = 19; Doc = "9148752";  Nombre = "HUMBERTO DE JESUS MULET HERNANDEZ";   Periodo = "1802"; Mora = 31249 },
    @{ Row = 20; Doc = "19896330"; Nombre = "CARLOS ALEXANDER FERNANDEZ CASTILLO"; Periodo = "1802"; Mora = 31249 },
    @{ Row = 21; Doc = "22800539"; Nombre = "DARLY YANETH PEREA TUIRAN";           Periodo = "1803"; Mora = 31249 },
    @{ Row = 22; Doc = "73099457"; Nombre = "RAFAEL SALAZAR MACHADO";             Periodo = "1803"; Mora = 31249 },
    @{ Row = 23; Doc = "9148752";  Nombre = "HUMBERTO DE JESUS MULET HERNANDEZ";   Periodo = "1803"; Mora = 31249 },
    @{ Row = 24; Doc = "19896330"; Nombre = "CARLOS ALEXANDER FERNANDEZ CASTILLO"; Periodo = "1803"; Mora = 31249 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
}
